$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shape = $s.Shapes.Item(3)
$seq = $s.TimeLine.MainSequence
$ids = 54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100
foreach ($id in $ids) {
    try {
        $eff = $seq.AddEffect($shape, $id, 0, 1)
        Write-Host $id "->" $eff.DisplayName
    } catch {
        Write-Host $id "-> ERROR" $_.Exception.Message
    }
}
